$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(469,44588,17779,3132,20911),
    @(470,44589,17796,3137,20933),
    @(471,44590,17817,3144,20961),
    @(472,44591,17830,3146,20976),
    @(473,44592,17850,3150,21000),
    @(474,44593,17877,3161,21038),
    @(475,44594,17896,3164,21060),
    @(476,44595,17921,3168,21089),
    @(477,44596,17938,3173,21111),
    @(478,44597,17958,3179,21137),
    @(479,44598,17973,3182,21155),
    @(480,44599,17995,3210,21205),
    @(481,44600,18014,3220,21234),
    @(482,44601,18040,3224,21264),
    @(483,44602,18064,3236,21300),
    @(484,44603,18081,3247,21328),
    @(485,44604,18095,3256,21351),
    @(486,44605,18105,3258,21363),
    @(487,44606,18123,3268,21391),
    @(488,44607,18145,3282,21427),
    @(489,44608,18179,3298,21477),
    @(490,44609,18201,3309,21510),
    @(491,44610,18225,3323,21548),
    @(492,44611,18240,3333,21573),
    @(493,44612,18252,3338,21590),
    @(494,44613,18279,3348,21627),
    @(495,44614,18314,3368,21682),
    @(496,44615,18338,3378,21716),
    @(497,44616,18375,3391,21766),
    @(498,44617,18413,3405,21818),
    @(499,44618,18453,3415,21868),
    @(500,44619,18485,3425,21910),
    @(501,44620,18530,3441,21971),
    @(502,44621,18567,3451,22018),
    @(503,44622,18611,3467,22078),
    @(504,44623,18631,3480,22111),
    @(505,44624,18663,3488,22151),
    @(506,44625,18687,3495,22182),
    @(507,44626,18704,3496,22200),
    @(508,44627,18750,3510,22260),
    @(509,44628,18786,3523,22309),
    @(510,44629,18817,3535,22352),
    @(511,44630,18847,3551,22398),
    @(512,44631,18881,3578,22459),
    @(513,44632,18918,3602,22520),
    @(514,44633,18950,3624,22574),
    @(515,44634,18979,3639,22618),
    @(516,44635,19009,3649,22658),
    @(517,44636,19040,3658,22698),
    @(518,44637,19073,3671,22744),
    @(519,44638,19093,3680,22773),
    @(520,44639,19110,3684,22794),
    @(521,44640,19133,3687,22820),
    @(522,44641,19161,3695,22856),
    @(523,44642,19185,3704,22889),
    @(524,44643,19209,3712,22921),
    @(525,44644,19229,3719,22948),
    @(526,44645,19251,3727,22978),
    @(527,44646,19270,3733,23003),
    @(528,44647,19292,3739,23031),
    @(529,44648,19311,3743,23054),
    @(530,44649,19334,3753,23087),
    @(531,44650,19352,3760,23112),
    @(532,44651,19368,3769,23137),
    @(533,44652,19392,3783,23175),
    @(534,44653,19417,3795,23212),
    @(535,44654,19440,3812,23252),
    @(536,44655,19462,3824,23286),
    @(537,44656,19482,3833,23315),
    @(538,44657,19500,3844,23344)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Scroll the view so the newly-added tail of the sheet is visible, then
# land the selection on the final data row/cell (matches the author's
# end-of-edit cursor position).
$excel.ActiveWindow.ScrollRow = 489
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A538").Select()
